$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new daily row (row 60) after the existing last row (59).
$newRow = 60

$ws.Cells.Item($newRow, 1).Value = 46009
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = 127
$ws.Cells.Item($newRow, 3).Value = 142
$ws.Cells.Item($newRow, 4).Value = 132
